$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("I3").Value = 4.7
$ws.Range("L3").Value = 4.85
$ws.Range("X3").Value = 7.2
$ws.Range("AH3").Value = 12
$ws.Range("AI3").Value = 26
$ws.Range("AJ3").Value = 15.5
$ws.Range("AK3").Value = 80
$ws.Range("AM3").Value = 60

# Row 4
$ws.Range("G4").Value = 1.24
$ws.Range("H4").Value = 5.3
$ws.Range("J4").Value = 1.62
$ws.Range("K4").Value = 2.55
$ws.Range("L4").Value = 8.75
$ws.Range("O4").Value = 1.53
$ws.Range("P4").Value = 2.2
$ws.Range("Q4").Value = 2.25
$ws.Range("U4").Value = 2.02
$ws.Range("V4").Value = 1.62
$ws.Range("W4").Value = 7.5
$ws.Range("X4").Value = 6.2
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 7.3
$ws.Range("AB4").Value = 30
$ws.Range("AC4").Value = 14
$ws.Range("AD4").Value = 11
$ws.Range("AE4").Value = 26
$ws.Range("AF4").Value = 120
$ws.Range("AH4").Value = 27
$ws.Range("AI4").Value = 80
$ws.Range("AM4").Value = 120

# Row 6
$ws.Range("G6").Value = 2.37
$ws.Range("I6").Value = 2.72
$ws.Range("J6").Value = 2.92
$ws.Range("L6").Value = 3.25
$ws.Range("N6").Value = 2.8
$ws.Range("T6").Value = 2.6
$ws.Range("W6").Value = 7.5
$ws.Range("X6").Value = 11
$ws.Range("Y6").Value = 9.5
$ws.Range("Z6").Value = 24
$ws.Range("AA6").Value = 20
$ws.Range("AD6").Value = 6.4
$ws.Range("AH6").Value = 8.25
$ws.Range("AI6").Value = 13
$ws.Range("AJ6").Value = 10.25
$ws.Range("AK6").Value = 30
$ws.Range("AL6").Value = 24
$ws.Range("AM6").Value = 35

# Row 10
$ws.Range("G10").Value = 2.3
$ws.Range("I10").Value = 3.2
$ws.Range("S10").Value = 1.44
$ws.Range("T10").Value = 2.63
$ws.Range("AC10").Value = 9
$ws.Range("AG10").Value = 201
$ws.Range("AO10").Value = 10

# Row 11
$ws.Range("H11").Value = 2.9
$ws.Range("O11").Value = 2.5
$ws.Range("P11").Value = 1.5
$ws.Range("S11").Value = 1.57
$ws.Range("T11").Value = 2.25
$ws.Range("X11").Value = 9.5
$ws.Range("AP11").Value = 1.88
$ws.Range("AQ11").Value = 1.98

# Row 12
$ws.Range("G12").Value = 3.6
$ws.Range("I12").Value = 2.15
$ws.Range("J12").Value = 4.33
$ws.Range("L12").Value = 3
$ws.Range("W12").Value = 8.5
$ws.Range("X12").Value = 17
$ws.Range("Y12").Value = 13
$ws.Range("Z12").Value = 41
$ws.Range("AA12").Value = 34
$ws.Range("AD12").Value = 6
$ws.Range("AE12").Value = 17
$ws.Range("AH12").Value = 6
$ws.Range("AI12").Value = 9
$ws.Range("AJ12").Value = 9.5
$ws.Range("AK12").Value = 19
$ws.Range("AL12").Value = 21
$ws.Range("AP12").Value = 1.83
$ws.Range("AQ12").Value = 2.03

# Row 16
$ws.Range("I16").Value = 3.15
$ws.Range("J16").Value = 3
$ws.Range("K16").Value = 1.93
$ws.Range("M16").Value = 1.5
$ws.Range("N16").Value = 2.42
$ws.Range("O16").Value = 2.42
$ws.Range("P16").Value = 1.5
$ws.Range("Q16").Value = 4.35
$ws.Range("R16").Value = 1.18
$ws.Range("S16").Value = 1.52
$ws.Range("T16").Value = 2.35
$ws.Range("U16").Value = 2.02
$ws.Range("V16").Value = 1.7
$ws.Range("W16").Value = 6.2
$ws.Range("X16").Value = 10.25
$ws.Range("Y16").Value = 9.5
$ws.Range("AA16").Value = 23
$ws.Range("AB16").Value = 40
$ws.Range("AC16").Value = 5.6
$ws.Range("AD16").Value = 5.8
$ws.Range("AE16").Value = 17
$ws.Range("AF16").Value = 110
$ws.Range("AH16").Value = 7.3
$ws.Range("AI16").Value = 14.5
$ws.Range("AJ16").Value = 11.75
$ws.Range("AK16").Value = 40
$ws.Range("AL16").Value = 35
$ws.Range("AM16").Value = 50
$ws.Range("AN16").Value = 1.11
$ws.Range("AO16").Value = 5.6
